$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -18.25145004873754
$ws.Cells.Item(2, 3).Value = -0.7398440805401197
$ws.Cells.Item(2, 4).Value = -18.25145004873754
$ws.Cells.Item(2, 5).Value = -18.25145004873754
$ws.Cells.Item(2, 6).Value = -18.25145004873754
$ws.Cells.Item(2, 7).Value = -18.25145004873754
$ws.Cells.Item(2, 8).Value = -18.25145004873754
$ws.Cells.Item(2, 9).Value = -18.25145004873754
$ws.Cells.Item(2, 10).Value = -18.25145004873754
$ws.Cells.Item(2, 11).Value = -18.25145004873754

$ws.Cells.Item(3, 2).Value = -18.25145004873754
$ws.Cells.Item(3, 3).Value = -18.25145004873754
$ws.Cells.Item(3, 4).Value = -18.25145004873754
$ws.Cells.Item(3, 5).Value = -18.25145004873754
$ws.Cells.Item(3, 6).Value = -18.25145004873754
$ws.Cells.Item(3, 7).Value = -18.25145004873754
$ws.Cells.Item(3, 8).Value = -18.25145004873754
$ws.Cells.Item(3, 9).Value = 0.721756822802923
$ws.Cells.Item(3, 10).Value = -18.25145004873754
$ws.Cells.Item(3, 11).Value = -18.25145004873754

$ws.Cells.Item(4, 2).Value = -18.25145004873754
$ws.Cells.Item(4, 3).Value = -0.6799566183519988
$ws.Cells.Item(4, 4).Value = -0.08824755012862201
$ws.Cells.Item(4, 5).Value = -18.25145004873754
$ws.Cells.Item(4, 6).Value = 4.021819969125364
$ws.Cells.Item(4, 7).Value = -18.25145004873754
$ws.Cells.Item(4, 8).Value = 2.884460223847127
$ws.Cells.Item(4, 9).Value = -18.25145004873754
$ws.Cells.Item(4, 10).Value = 4.321923702861627
$ws.Cells.Item(4, 11).Value = -18.25145004873754

$ws.Cells.Item(5, 2).Value = -18.25145004873754
$ws.Cells.Item(5, 3).Value = -0.006375559567528655
$ws.Cells.Item(5, 4).Value = -18.25145004873754
$ws.Cells.Item(5, 5).Value = -18.25145004873754
$ws.Cells.Item(5, 6).Value = -18.25145004873754
$ws.Cells.Item(5, 7).Value = 3.87783328268995
$ws.Cells.Item(5, 8).Value = -18.25145004873754
$ws.Cells.Item(5, 9).Value = -18.25145004873754
$ws.Cells.Item(5, 10).Value = -18.25145004873754
$ws.Cells.Item(5, 11).Value = -18.25145004873754

$ws.Cells.Item(6, 2).Value = -18.25145004873754
$ws.Cells.Item(6, 3).Value = -18.25145004873754
$ws.Cells.Item(6, 4).Value = -18.25145004873754
$ws.Cells.Item(6, 5).Value = -18.25145004873754
$ws.Cells.Item(6, 6).Value = -18.25145004873754
$ws.Cells.Item(6, 7).Value = -18.25145004873754
$ws.Cells.Item(6, 8).Value = -18.25145004873754
$ws.Cells.Item(6, 9).Value = -18.25145004873754
$ws.Cells.Item(6, 10).Value = -18.25145004873754
$ws.Cells.Item(6, 11).Value = -18.25145004873754

$ws.Cells.Item(7, 2).Value = 3.304556823647343
$ws.Cells.Item(7, 3).Value = -18.25145004873754
$ws.Cells.Item(7, 4).Value = -18.25145004873754
$ws.Cells.Item(7, 5).Value = -18.25145004873754
$ws.Cells.Item(7, 6).Value = -18.25145004873754
$ws.Cells.Item(7, 7).Value = -18.25145004873754
$ws.Cells.Item(7, 8).Value = -18.25145004873754
$ws.Cells.Item(7, 9).Value = -18.25145004873754
$ws.Cells.Item(7, 10).Value = -18.25145004873754
$ws.Cells.Item(7, 11).Value = -18.25145004873754

$ws.Cells.Item(8, 2).Value = -18.25145004873754
$ws.Cells.Item(8, 3).Value = -18.25145004873754
$ws.Cells.Item(8, 4).Value = -18.25145004873754
$ws.Cells.Item(8, 5).Value = 1.219075774669758
$ws.Cells.Item(8, 6).Value = -18.25145004873754
$ws.Cells.Item(8, 7).Value = -18.25145004873754
$ws.Cells.Item(8, 8).Value = -18.25145004873754
$ws.Cells.Item(8, 9).Value = -18.25145004873754
$ws.Cells.Item(8, 10).Value = -18.25145004873754
$ws.Cells.Item(8, 11).Value = -18.25145004873754

$ws.Cells.Item(9, 2).Value = 3.339084464093891
$ws.Cells.Item(9, 3).Value = -18.25145004873754
$ws.Cells.Item(9, 4).Value = -18.25145004873754
$ws.Cells.Item(9, 5).Value = -18.25145004873754
$ws.Cells.Item(9, 6).Value = -18.25145004873754
$ws.Cells.Item(9, 7).Value = -18.25145004873754
$ws.Cells.Item(9, 8).Value = -18.25145004873754
$ws.Cells.Item(9, 9).Value = -18.25145004873754
$ws.Cells.Item(9, 10).Value = -18.25145004873754
$ws.Cells.Item(9, 11).Value = -18.25145004873754

$ws.Cells.Item(10, 2).Value = -18.25145004873754
$ws.Cells.Item(10, 3).Value = -18.25145004873754
$ws.Cells.Item(10, 4).Value = -18.25145004873754
$ws.Cells.Item(10, 5).Value = -18.25145004873754
$ws.Cells.Item(10, 6).Value = -18.25145004873754
$ws.Cells.Item(10, 7).Value = -18.25145004873754
$ws.Cells.Item(10, 8).Value = -18.25145004873754
$ws.Cells.Item(10, 9).Value = -0.3669158006929812
$ws.Cells.Item(10, 10).Value = -18.25145004873754
$ws.Cells.Item(10, 11).Value = 1.945148476037946

$ws.Cells.Item(11, 2).Value = -18.25145004873754
$ws.Cells.Item(11, 3).Value = -18.25145004873754
$ws.Cells.Item(11, 4).Value = -18.25145004873754
$ws.Cells.Item(11, 5).Value = 2.622682569334144
$ws.Cells.Item(11, 6).Value = -18.25145004873754
$ws.Cells.Item(11, 7).Value = 1.387786425389846
$ws.Cells.Item(11, 8).Value = -18.25145004873754
$ws.Cells.Item(11, 9).Value = -18.25145004873754
$ws.Cells.Item(11, 10).Value = -18.25145004873754
$ws.Cells.Item(11, 11).Value = 0.716308137348333

$ws.Cells.Item(12, 2).Value = -18.25145004873754
$ws.Cells.Item(12, 3).Value = -18.25145004873754
$ws.Cells.Item(12, 4).Value = -18.25145004873754
$ws.Cells.Item(12, 5).Value = -18.25145004873754
$ws.Cells.Item(12, 6).Value = -18.25145004873754
$ws.Cells.Item(12, 7).Value = -18.25145004873754
$ws.Cells.Item(12, 8).Value = -18.25145004873754
$ws.Cells.Item(12, 9).Value = -18.25145004873754
$ws.Cells.Item(12, 10).Value = -18.25145004873754
$ws.Cells.Item(12, 11).Value = -18.25145004873754

$ws.Cells.Item(13, 2).Value = -18.25145004873754
$ws.Cells.Item(13, 3).Value = -18.25145004873754
$ws.Cells.Item(13, 4).Value = -18.25145004873754
$ws.Cells.Item(13, 5).Value = 2.032810388639164
$ws.Cells.Item(13, 6).Value = -18.25145004873754
$ws.Cells.Item(13, 7).Value = -18.25145004873754
$ws.Cells.Item(13, 8).Value = -18.25145004873754
$ws.Cells.Item(13, 9).Value = -18.25145004873754
$ws.Cells.Item(13, 10).Value = -18.25145004873754
$ws.Cells.Item(13, 11).Value = 3.417606708504008

$ws.Cells.Item(14, 2).Value = -18.25145004873754
$ws.Cells.Item(14, 3).Value = -18.25145004873754
$ws.Cells.Item(14, 4).Value = 0.7231577942262535
$ws.Cells.Item(14, 5).Value = -18.25145004873754
$ws.Cells.Item(14, 6).Value = -18.25145004873754
$ws.Cells.Item(14, 7).Value = -18.25145004873754
$ws.Cells.Item(14, 8).Value = -18.25145004873754
$ws.Cells.Item(14, 9).Value = -18.25145004873754
$ws.Cells.Item(14, 10).Value = -18.25145004873754
$ws.Cells.Item(14, 11).Value = 1.112267653815624

$ws.Cells.Item(15, 2).Value = -18.25145004873754
$ws.Cells.Item(15, 3).Value = -18.25145004873754
$ws.Cells.Item(15, 4).Value = -0.5712759841182389
$ws.Cells.Item(15, 5).Value = -18.25145004873754
$ws.Cells.Item(15, 6).Value = -18.25145004873754
$ws.Cells.Item(15, 7).Value = -18.25145004873754
$ws.Cells.Item(15, 8).Value = -18.25145004873754
$ws.Cells.Item(15, 9).Value = -18.25145004873754
$ws.Cells.Item(15, 10).Value = -18.25145004873754
$ws.Cells.Item(15, 11).Value = -18.25145004873754

$ws.Cells.Item(16, 2).Value = -18.25145004873754
$ws.Cells.Item(16, 3).Value = -18.25145004873754
$ws.Cells.Item(16, 4).Value = -18.25145004873754
$ws.Cells.Item(16, 5).Value = -18.25145004873754
$ws.Cells.Item(16, 6).Value = -18.25145004873754
$ws.Cells.Item(16, 7).Value = -18.25145004873754
$ws.Cells.Item(16, 8).Value = -18.25145004873754
$ws.Cells.Item(16, 9).Value = -18.25145004873754
$ws.Cells.Item(16, 10).Value = -18.25145004873754
$ws.Cells.Item(16, 11).Value = -18.25145004873754

$ws.Cells.Item(17, 2).Value = -18.25145004873754
$ws.Cells.Item(17, 3).Value = 0.1910439923710072
$ws.Cells.Item(17, 4).Value = -0.4146784771341909
$ws.Cells.Item(17, 5).Value = -18.25145004873754
$ws.Cells.Item(17, 6).Value = -18.25145004873754
$ws.Cells.Item(17, 7).Value = -18.25145004873754
$ws.Cells.Item(17, 8).Value = 1.605268777151363
$ws.Cells.Item(17, 9).Value = -1.122156968385704
$ws.Cells.Item(17, 10).Value = -18.25145004873754
$ws.Cells.Item(17, 11).Value = -18.25145004873754

$ws.Cells.Item(18, 2).Value = -18.25145004873754
$ws.Cells.Item(18, 3).Value = -18.25145004873754
$ws.Cells.Item(18, 4).Value = -18.25145004873754
$ws.Cells.Item(18, 5).Value = -18.25145004873754
$ws.Cells.Item(18, 6).Value = -18.25145004873754
$ws.Cells.Item(18, 7).Value = -18.25145004873754
$ws.Cells.Item(18, 8).Value = 2.077856874545243
$ws.Cells.Item(18, 9).Value = -1.447962150526292
$ws.Cells.Item(18, 10).Value = -18.25145004873754
$ws.Cells.Item(18, 11).Value = -18.25145004873754

$ws.Cells.Item(19, 2).Value = -18.25145004873754
$ws.Cells.Item(19, 3).Value = -18.25145004873754
$ws.Cells.Item(19, 4).Value = 2.840910305776042
$ws.Cells.Item(19, 5).Value = -18.25145004873754
$ws.Cells.Item(19, 6).Value = -18.25145004873754
$ws.Cells.Item(19, 7).Value = -18.25145004873754
$ws.Cells.Item(19, 8).Value = 1.492610174198279
$ws.Cells.Item(19, 9).Value = -0.3351915175671321
$ws.Cells.Item(19, 10).Value = -18.25145004873754
$ws.Cells.Item(19, 11).Value = -18.25145004873754

$ws.Cells.Item(20, 2).Value = -18.25145004873754
$ws.Cells.Item(20, 3).Value = 3.427248511028375
$ws.Cells.Item(20, 4).Value = 3.140872091271929
$ws.Cells.Item(20, 5).Value = -18.25145004873754
$ws.Cells.Item(20, 6).Value = 1.90924002984275
$ws.Cells.Item(20, 7).Value = -18.25145004873754
$ws.Cells.Item(20, 8).Value = 0.2353103541553732
$ws.Cells.Item(20, 9).Value = 3.996090723461194
$ws.Cells.Item(20, 10).Value = -18.25145004873754
$ws.Cells.Item(20, 11).Value = 0.7300308303048963

$ws.Cells.Item(21, 2).Value = -18.25145004873754
$ws.Cells.Item(21, 3).Value = 2.556418174776356
$ws.Cells.Item(21, 4).Value = -18.25145004873754
$ws.Cells.Item(21, 5).Value = 2.891617346818376
$ws.Cells.Item(21, 6).Value = -18.25145004873754
$ws.Cells.Item(21, 7).Value = 1.423454244670587
$ws.Cells.Item(21, 8).Value = 0.4438224289678785
$ws.Cells.Item(21, 9).Value = -18.25145004873754
$ws.Cells.Item(21, 10).Value = -18.25145004873754
$ws.Cells.Item(21, 11).Value = -18.25145004873754

